$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 245; existing rows 245-296 shift down to 246-297.
$ws.Range("A245").EntireRow.Insert()

# Populate the newly inserted row 245 with the new price-observation record.
# (Same market/category/variety/quality/unit/origin as the record that used
# to sit at row 245, but a new date and new min/max/weighted-avg/kg prices.)
$ws.Range("A245").Value = 7
$ws.Range("B245").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C245").Value = "Ñuble"
$ws.Range("D245").Value = 44476
$ws.Range("E245").Value = 16
$ws.Range("F245").Value = 100112020
$ws.Range("G245").Value = "Tomate"
$ws.Range("H245").Value = "Larga vida"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 160
$ws.Range("K245").Value = 7500
$ws.Range("L245").Value = 8000
$ws.Range("M245").Value = 7750
$ws.Range("N245").Value = "`$/caja 10 kilos"
$ws.Range("O245").Value = "Región de Arica y Parinacota"
$ws.Range("P245").Value = 775
$ws.Range("Q245").Value = 10
$ws.Range("R245").Value = "Hortaliza"
